$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9: mark E9 (Status) as Completed, and clear the stray orange-fill
# formatting on F9 (empty cell that previously only carried a fill style).
$ws.Range("E9").Value = "Completed"
$ws.Range("F9").Clear()

# Row 14: reassign Owner (D14) to Stefan and mark Status (E14) Completed.
$ws.Range("D14").Value = "Stefan"
$ws.Range("E14").Value = "Completed"

# Row 21: mark E21 (Status) as Completed.
$ws.Range("E21").Value = "Completed"

# Update the view: scroll so row 10 is the top-left visible row, and move
# the active selection to G21.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G21").Select()
